$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "展览" (Exhibition) - update "想去人数" (F) and "最低票价" (G) values
# ---------------------------------------------------------------------------
$wsExpo = $wb.Worksheets("展览")

$wsExpo.Range("F4").Value = 5953
$wsExpo.Range("F5").Value = 5953
$wsExpo.Range("F6").Value = 71
$wsExpo.Range("F13").Value = 298
$wsExpo.Range("F14").Value = 4396
$wsExpo.Range("F15").Value = 4396
$wsExpo.Range("G16").Value = "不可售"
$wsExpo.Range("G20").Value = "不可售"
$wsExpo.Range("F22").Value = 6798
$wsExpo.Range("F23").Value = 6798
$wsExpo.Range("G25").Value = "不可售"
$wsExpo.Range("F26").Value = 466
$wsExpo.Range("F27").Value = 1259
$wsExpo.Range("F32").Value = 6015
$wsExpo.Range("F37").Value = 430
$wsExpo.Range("F38").Value = 4855
$wsExpo.Range("F50").Value = 2075
$wsExpo.Range("F51").Value = 14
$wsExpo.Range("F52").Value = 1029

# ---------------------------------------------------------------------------
# Sheet "演出" (Performance) - update F5, then insert a new event row (row 8)
# pushing the previous rows 8-13 down to 9-14.
# ---------------------------------------------------------------------------
$wsShow = $wb.Worksheets("演出")

$wsShow.Range("F5").Value = 32

# Insert a brand-new blank row at position 8 (shifts rows 8-13 -> 9-14)
$wsShow.Rows.Item(8).Insert()

# Give the new row's index cell (column A) the same look as the other rows
$idxCell = $wsShow.Cells.Item(8, 1)
$idxCell.Font.Bold = $true
$idxCell.HorizontalAlignment = -4108
$idxCell.VerticalAlignment = -4160
$idxCell.Borders.LineStyle = 1
$idxCell.Value = 7

# Column B holds a plain text date string ("2024-04-21") - force text so Excel
# doesn't silently reinterpret it as a date serial number.
$dateCell = $wsShow.Cells.Item(8, 2)
$dateCell.NumberFormat = "@"
$dateCell.Value = "2024-04-21"
$dateCell.Style = "Normal"

$wsShow.Cells.Item(8, 3).Value = "北京·次元音浪 Million Live—超新星派对"
$wsShow.Cells.Item(8, 4).Value = "学清路38号金码大厦B座 北京想象空间"
$wsShow.Cells.Item(8, 5).Value = "2024.04.21 13:00-04.21 16:00"
$wsShow.Cells.Item(8, 6).Value = 1
$wsShow.Cells.Item(8, 7).Value = 88
$wsShow.Cells.Item(8, 8).Value = "https://show.bilibili.com/platform/detail.html?id=83202"
$wsShow.Cells.Item(8, 9).Value = "//i0.hdslb.com/bfs/openplatform/202403/OfpyKpSQ1711013512280.png"

# ---------------------------------------------------------------------------
# Sheet "全部类型" (All types) - update "想去人数" (F) and "最低票价" (G) values
# ---------------------------------------------------------------------------
$wsAll = $wb.Worksheets("全部类型")

$wsAll.Range("F4").Value = 5953
$wsAll.Range("F5").Value = 5953
$wsAll.Range("F6").Value = 71
$wsAll.Range("F12").Value = 197
$wsAll.Range("F13").Value = 298
$wsAll.Range("F14").Value = 4396
$wsAll.Range("F15").Value = 4396
$wsAll.Range("G16").Value = "不可售"
$wsAll.Range("G20").Value = "不可售"
$wsAll.Range("F22").Value = 6798
$wsAll.Range("F23").Value = 6798
$wsAll.Range("G25").Value = "不可售"
$wsAll.Range("F26").Value = 466
$wsAll.Range("F27").Value = 1259
$wsAll.Range("F34").Value = 6015
$wsAll.Range("F39").Value = 430
$wsAll.Range("F40").Value = 4856
$wsAll.Range("F51").Value = 14
